$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 5) of data to the table, mirroring the existing
# login/contact rows above it.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "el"
$ws.Range("C5").Value = "T."
$ws.Range("D5").Value = "Gil"
$ws.Range("E5").Value = "e@gmail.com"

# Add a mailto hyperlink on the new email cell, matching the style used
# by the other email cells (E2:E4).
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:e@gmail.com") | Out-Null
$ws.Range("E5").Style = "Hyperlink"

# Move the active selection to the new last cell, matching the saved
# workbook's last selection state.
$ws.Range("E5").Select() | Out-Null
